# Update gh-pages to output generated at 456a3b4
#
# Sheet "展览" (index 1): update a handful of "want to go" counters (col F)
# and append a new row (row 10) for the 九江 event.
#
# Sheet "全部类型" (index 4): same counter updates and the same new row,
# appended as row 11 (this sheet already had the extra row that "展览"
# gained, plus one more event in between).

function Set-NewEventRow {
    param($ws, $rowIndex, $indexValue)

    # Copy formatting (incl. the bold/centered/bordered style) from the
    # previous row's column-A cell so the new row matches existing rows.
    $ws.Cells.Item($rowIndex - 1, 1).Copy($ws.Cells.Item($rowIndex, 1))
    $ws.Cells.Item($rowIndex, 1).Value = $indexValue

    # Keep date-looking text as plain text, not auto-converted to a date.
    $ws.Cells.Item($rowIndex, 2).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 2).Value = "2025-02-14"

    $ws.Cells.Item($rowIndex, 3).Value = "九江·第二届异次元动漫嘉年华"
    $ws.Cells.Item($rowIndex, 4).Value = "长虹西大道兴城广场99号 九江半岛宾馆"

    $ws.Cells.Item($rowIndex, 5).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 5).Value = "2025.02.14 09:30-02.14 17:30"

    $ws.Cells.Item($rowIndex, 6).Value = 2
    $ws.Cells.Item($rowIndex, 7).Value = 39.8
    $ws.Cells.Item($rowIndex, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93138"
    $ws.Cells.Item($rowIndex, 9).Value = "//i1.hdslb.com/bfs/openplatform/202409/YBlAWRDD1727019019550.jpeg"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 6).Value = 717   # 712 -> 717
$ws1.Cells.Item(3, 6).Value = 36    # 34  -> 36
$ws1.Cells.Item(4, 6).Value = 241   # 240 -> 241
$ws1.Cells.Item(5, 6).Value = 2557  # 2502 -> 2557
$ws1.Cells.Item(7, 6).Value = 3626  # 3591 -> 3626
$ws1.Cells.Item(9, 6).Value = 918   # 910 -> 918

Set-NewEventRow $ws1 10 9

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2, 6).Value = 717    # 712 -> 717
$ws4.Cells.Item(3, 6).Value = 36     # 34  -> 36
$ws4.Cells.Item(5, 6).Value = 241    # 240 -> 241
$ws4.Cells.Item(6, 6).Value = 2557   # 2502 -> 2557
$ws4.Cells.Item(8, 6).Value = 3626   # 3591 -> 3626
$ws4.Cells.Item(10, 6).Value = 918   # 910 -> 918

Set-NewEventRow $ws4 11 10
